$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 174 updates
$ws.Range("B174").Value = 7057
$ws.Range("D174").Value = 57283
$ws.Range("E174").Value = 121623
$ws.Range("K174").Value = 14656
$ws.Range("L174").Value = 58982

# Row 175 updates
$ws.Range("B175").Value = 6760
$ws.Range("D175").Value = 56655
$ws.Range("E175").Value = 120070
$ws.Range("H175").Value = 16011
$ws.Range("I175").Value = 7277
$ws.Range("K175").Value = 14889
$ws.Range("L175").Value = 59776

# Row 176 updates
$ws.Range("B176").Value = 7225
$ws.Range("D176").Value = 56858
$ws.Range("E176").Value = 120941
$ws.Range("G176").Value = 5524
$ws.Range("H176").Value = 18530
$ws.Range("I176").Value = 7572
$ws.Range("K176").Value = 14625
$ws.Range("L176").Value = 61039

# Row 177: previously only had A177, C177, D177 populated. Fill in the rest.
$ws.Range("B177").Value = 6183
$ws.Range("D177").Value = 60503
$ws.Range("E177").Value = 127189
$ws.Range("F177").Value = 1017
$ws.Range("G177").Value = 8663
$ws.Range("H177").Value = 18981
$ws.Range("I177").Value = 8792
$ws.Range("J177").Value = 13420
$ws.Range("K177").Value = 14578
$ws.Range("L177").Value = 61739
